$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.115.06"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.642.08"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.143"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000190"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "3.127.54"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "68.047.01"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "2.645.29"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "362.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").Value = "2.778.63"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "553.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "0.0₆0338"
$ws.Range("E42").Value = "  +5.60%  "
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").Value = "  +0.60%  "
